# Add two new columns (I0, IF) to the right of the existing "IP" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: copy the formatting of H1 ("IP" header) onto I1:J1, then set text ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-36: "I0" and "IF" values per row ---
$i0 = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 3; 7 = 2; 8 = 1; 9 = 4; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 1; 27 = 1; 28 = 1;
    29 = 1; 30 = 1; 31 = 1; 32 = 1; 33 = 1; 34 = 1; 35 = 1; 36 = 1
}
$if = @{
    2 = 6; 3 = 5; 4 = 4; 5 = 2; 6 = 6; 7 = 5; 8 = 5; 9 = 7; 10 = 4;
    11 = 5; 12 = 5; 13 = 4; 14 = 6; 15 = 6; 16 = 5; 17 = 5; 18 = 6; 19 = 5;
    20 = 5; 21 = 5; 22 = 6; 23 = 5; 24 = 6; 25 = 5; 26 = 4; 27 = 5; 28 = 6;
    29 = 5; 30 = 6; 31 = 6; 32 = 6; 33 = 3; 34 = 4; 35 = 3; 36 = 2
}

for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 9).Value = $i0[$r]
    $ws.Cells.Item($r, 10).Value = $if[$r]
}

Write-Output "I0 and IF columns added"
